$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 5) with the same shape/style as the existing
# rows, matching the weekly "Fruta / hortaliza" price update for Vega
# Monumental Concepción - Tuna.
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"

$ws.Range("D5").Value = 44516
$ws.Range("D5").NumberFormat = $ws.Range("D4").NumberFormat

$ws.Range("E5").Value = 8
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107011
$ws.Range("J5").Value = "Tuna"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 33000
$ws.Range("O5").Value = 34000
$ws.Range("P5").Value = 33500
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Provincia de Melipilla"
$ws.Range("S5").Value = 1861
$ws.Range("T5").Value = 18
